$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# A new column ("caiman") is inserted before the existing "AWS"
# column (K). We replicate that by first relocating the existing
# K-column ("AWS") values to column L (preserving the "AWS" shared
# string reference so it keeps its position in the shared string
# table), then populating the new K column with the "todo"/"caiman"
# values from the diff. Formatting (style index 1, the sheet's usual
# centered style) is (re)applied explicitly by pasting the format
# from the neighboring I-column cell on the same row, which is
# reliably styled with index 1 throughout the sheet.
# -----------------------------------------------------------------

# 1) Move the "AWS" header text and its boolean flags from K to L.
$ws.Range("L1").Value = "AWS"
$boolRows = 3,4,5,6,7,9,10,11,12,13,14,15
foreach ($r in $boolRows) {
    $ws.Range("L$r").Value = $true
}

# Re-apply the sheet's normal centered style (index 1) to every cell
# we just touched in column L, using the always-style-1 column I cell
# on the same row as the format source.
$lRows = @(1) + $boolRows
foreach ($r in $lRows) {
    $ws.Range("I$r").Copy()
    $ws.Range("L$r").PasteSpecial(-4122)
}

# Clear out the old K-column cells that had no replacement value so
# that they fully disappear (no stray empty cell left behind).
$clearKRows = 3,4,6,7,9,10,11,12,13,14,15
foreach ($r in $clearKRows) {
    $ws.Range("K$r").Clear()
}

# 2) Populate the new "caiman" column (K) with its header and values.
$newK = @{1="caiman"; 2="todo"; 5="todo"; 8="todo"; 16="caiman"; 21="todo"; 22="todo"; 27="todo"; 28="todo"}
foreach ($r in $newK.Keys) {
    $ws.Range("K$r").Value = $newK[$r]
}
foreach ($r in $newK.Keys) {
    $ws.Range("I$r").Copy()
    $ws.Range("K$r").PasteSpecial(-4122)
}

$ws.Application.CutCopyMode = $false

# 3) Update the selection to match the saved view state.
$ws.Range("F15").Select()
